# Rewrite of example 11 unit tests
# Updates the "LookupValue UT Posting" / "LookupValue UT Inheritance" scenario
# descriptions from the old event-subscriber based naming to the new
# procedure-based naming, and adjusts the wrapped-text row heights to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ATDD Scenarios")

# --- Scenario (column F) and Given/When/Then description (column H) text ---

$ws.Range("F3").Value2  = "Check failure CheckLookupvalueExistsOnSalesHeader Sales Posting"
$ws.Range("H5").Value2  = "Trigger CheckLookupvalueExistsOnSalesHeader Sales Posting"
$ws.Range("H9").Value2  = "Trigger CheckLookupvalueExistsOnSalesHeader Sales Posting"

$ws.Range("F11").Value2 = "Check failure CheckLookupvalueExistsOnSalesHeader Whse. Posting"
$ws.Range("H13").Value2 = "Trigger CheckLookupvalueExistsOnSalesHeader Whse. Posting"

$ws.Range("F15").Value2 = "Check success CheckLookupvalueExistsOnSalesHeader Whse. Posting"
$ws.Range("H17").Value2 = "Trigger CheckLookupvalueExistsOnSalesHeader Whse. Posting"

$ws.Range("F20").Value2 = "Check InheritLookupValueFromCustomer"
$ws.Range("H23").Value2 = "Trigger InheritLookupValueFromCustomer"

$ws.Range("F25").Value2 = "Check ApplyLookupValueFromCustomerTemplate from Contact"
$ws.Range("H28").Value2 = "Trigger ApplyLookupValueFromCustomerTemplate"

$ws.Range("F30").Value2 = "Check ApplyLookupValueFromCustomerTemplate"
$ws.Range("H33").Value2 = "Trigger ApplyLookupValueFromCustomerTemplate"

$ws.Range("F36").Value2 = "Check InheritLookupValueFromSalesHeader"
$ws.Range("H39").Value2 = "Trigger InheritLookupValueFromSalesHeader"

# --- Row heights (wrapped text reflows because the new strings are a
#     different length than the old ones) ---

$ws.Rows.Item(3).RowHeight  = 45.75
$ws.Rows.Item(5).RowHeight  = 30
$ws.Rows.Item(7).RowHeight  = 45
$ws.Rows.Item(9).RowHeight  = 30
$ws.Rows.Item(11).RowHeight = 45
$ws.Rows.Item(13).RowHeight = 30
$ws.Rows.Item(15).RowHeight = 45
$ws.Rows.Item(17).RowHeight = 30
$ws.Rows.Item(20).RowHeight = 30.75
$ws.Rows.Item(23).RowHeight = 16.5
$ws.Rows.Item(25).RowHeight = 30
$ws.Rows.Item(28).RowHeight = 30
$ws.Rows.Item(30).RowHeight = 30
$ws.Rows.Item(36).RowHeight = 30.75

# --- Reset the active selection back to the top-left cell ---

$ws.Range("A1").Select()
